# "incorporated income groups into the analysis"
#
# The sheet gets one new row inserted above the header (with a blank
# spacer row under it) that documents the data source, and the last
# bucket's open-ended upper bound gets flagged (yellow highlight +
# explanatory comment) since it was an assumed cutoff rather than a
# value that came from the original distribution.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows at the very top of the sheet. Row 1 will hold the
# new source note, row 2 stays empty as a spacer before the (now shifted)
# header row, and every existing row/formula/merge below moves down by 2.
$ws.Range("A1:A2").EntireRow.Insert()

# Document where the digitized numbers came from.
$ws.Range("A1").Value = "Based on digitized data from this plot: https://upload.wikimedia.org/wikipedia/commons/c/cb/Distribution_of_Annual_Household_Income_in_the_United_States_2015.svg "

# The open-ended top bucket's upper bound (250000-1000000, now row 45) was
# assumed rather than taken from the source distribution - call it out.
$ws.Range("B45").Interior.Color = 65535
$comment = $ws.Range("B45").AddComment("Steffen Coenen:" + [char]10 + "assumed (cutoff wasn't given in the distribution)")

# Leave the selection where the author left it.
$ws.Range("A2").Select() | Out-Null
